$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New receipt row (row 10). Columns: B=Tarih C=Belge Adi D=Belge No
# E=Vergi Dairesi F=Vergi No G=KDV H=Tutar I=Toplam
# The numeric-looking values (D, F, G, H, I) must stay stored as TEXT
# (shared strings), exactly as every other cell in this table is stored
# as text. A leading apostrophe forces Excel to keep the literal text
# instead of auto-coercing it to a Number.
$ws.Range("B10").Value = "30/07/2022"
$ws.Range("C10").Value = "E ROL OFISI AS . ATAŞEHIR ŞB . MAPARK OTOPARK VE ULASIM HIZMETLERI SAN . VE TIC A.S NO : 460 "
$ws.Range("D10").Value = "'000900050006016521000192"
$ws.Range("E10").Value = "TARTH : 24/07/2022 FIS NO : 0165 2 AD x 335,00 TL / AD YENI RAKI 70 CL MIGROS PLASTIK POSET TOPKDV TOPLAM MIGROS TICARET A.S. MERKEZ SELCUK 12 IR MIGROS SATIS MGZSI ATATURK MH . KUSILAY CD . NO : 1 SELCUK / IZMIR TEL : 0 ( 85 ) 229 8094 MF FO 00053324 BUYUK MUKELLEFLER V.D.6220529513 MERKEZ ADRESI . ATATURK MAH . TURGUT OZAL BULV . NO : 7 ATASE - 13 / 1STANBUL 538121 ******* 3010 4521824 ****** 3017 ARA TOPLAM % 18 % 18 SAAT : 11 : 26 0.25 * 670 25 * 102,24 * 670,25 CRTAK POS ORTAK POS * 334,00 * 336.25 KOV TUTARI KDV Lİ TOPLAM S ASYON K NO : 0010 Z NO : 1441 MF DI 12110058 SEDIR RESTAURANT CEVAHIR ILHAN ŞİRİNCE MH . 204 SK . NO : 2 SELÇUK İZMİR 37993742150 AFİYET OLSUN TEŞEKKÜR EDERİZ 28-07-2022 SAAT : 22:24 YEMEK KDV TOP K.KARTI % 8 FIS NO : 10 TERMINAL 01922484 ( E ) ISYERI NO 0001500118 22:25 SATIS *********** 3010 LERKE ALTIPARMAK HID : AU000000041010 460,00 * 4,44 * 60,00 +60,00 Tutar : 60.00 TL KAKSILIGI MAL VEYA HİZMET ALINDI LUTFEN BU BELGEYI SAKLAYINIZ MASTERCARD KASA : KASIYER : E.BOZKUS http://www.mapark.com.tr Z No : 1093 TEŞEKKÜR EDERİZ CEK : 0053866 EKU No : 0001 MF SC 0000007733 : 30.07.2022 TARIH SAAT : 09:48 FIS NO : 0019 2 X 12,50 SIGARA BOREGI 400 G ALIŞVERİŞ POŞETİ TRAKYAESKİ KAŞAR300G ERITME PEYNİRİ 500 G L.SUZ SUT PINAR 1L 0,335 kg X 12,90 YEŞİL BİBER 0,890 kg X 12,90 DOMATES KARADENİZ EKMEĞİ SARIMSAK 250 G ODUN EKMEĞİ 5 X 4,45 SIMIT TOPKDV TOPLAM 1 ) BIM BIRLEŞIK MAĞAZALAR A.Ş Ağva Merkez Mah . Şile Cad . No : 171/1 Şile / İSTANBUL BUYUK MÜKELLEFLER "
$ws.Range("F10").Value = "'729001504362205295131750051846"
$ws.Range("G10").Value = "'3.85"
$ws.Range("H10").Value = "'2996.25"
$ws.Range("I10").Value = "'3000.10"
